# Edit script: add columns I (I0) and J (IF) with header + data values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (match style of existing header cells, e.g. H1: bold, thin border, centered/top aligned)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data values for rows 2-55, columns I (index 9) and J (index 10)
$data = @(
    @(3, 4),
    @(9, 9),
    @(7, 8),
    @(2, 3),
    @(5, 7),
    @(8, 8),
    @(6, 7),
    @(5, 5),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(6, 8),
    @(8, 8),
    @(8, 8),
    @(3, 5),
    @(5, 6),
    @(6, 6),
    @(9, 9),
    @(4, 6),
    @(5, 6),
    @(3, 6),
    @(8, 9),
    @(4, 5),
    @(5, 8),
    @(9, 9),
    @(7, 8),
    @(4, 6),
    @(6, 9),
    @(8, 9),
    @(5, 6),
    @(8, 9),
    @(3, 5),
    @(2, 5),
    @(11, 12),
    @(7, 7),
    @(6, 7),
    @(6, 8),
    @(6, 7),
    @(6, 7),
    @(4, 6),
    @(5, 6),
    @(7, 7),
    @(5, 7),
    @(7, 7),
    @(2, 3),
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(4, 5),
    @(2, 3),
    @(9, 9)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $data[$r][0]
    $ws.Cells.Item($rowNum, 10).Value = $data[$r][1]
}
